# Edit: append a follow-up sentence to the "Normally, you want to test
# hypothesis..." paragraph on the two "MDS Quality measure: Statistical
# descriptive" slides (slide 9 and slide 10), and grow their text box to
# match the extra line of wrapped text (spAutoFit-style resize).

$p = $ppt.ActivePresentation

# Curly apostrophe (U+2019) used by the existing "won't" in the deck - build
# the search string with it so it matches the existing single run exactly
# (and so we don't accidentally downgrade it to a straight quote).
$apos = [char]0x2019
$targetPhrase = "Normally, you want to test hypothesis to see if difference is statistically significant using Central Limit Theorem or any parametric tests. But this analysis shows that this method won" + $apos + "t work."
$suffix = " This analysis suggests me to use permutation test instead."

# New target height for the text box, in points (EMU 3970318 / 12700 EMU-per-point).
$newHeightPts = 3970318 / 12700

foreach ($slideIdx in @(9, 10)) {
    $slide = $p.Slides.Item($slideIdx)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)

        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $fullText = $tr.Text

            if ($fullText.Contains("Normally, you want")) {
                # Isolate the run of characters starting at "Normally, you want"
                # through the end of the text box, and rewrite just that span so
                # the paragraph stays a single run (matching the target markup).
                $startIdx = $fullText.IndexOf("Normally, you want")
                $totalLen = $fullText.Length
                $span = $tr.Characters($startIdx + 1, $totalLen - $startIdx)
                $span.Text = $targetPhrase + $suffix

                # Grow the shape to fit the now-longer paragraph.
                $shp.Height = $newHeightPts
            }
        }
    }
}
